$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C rows 2-496 all hold the "Förändrad" date and need to move
# from 45177 (2023-09-12) to 45178 (2023-09-13).
$ws.Range("C2:C496").Value = 45178
